$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 5.423951666666667
$ws.Range("H2").Value = 16.271855
$ws.Range("I2").Value = 0.4774188439413272
$ws.Range("J2").Value = 0.4774188439413271
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.707786666666667
$ws.Range("N2").Value = 5.12336
$ws.Range("O2").Value = 0.8764025646701329
$ws.Range("P2").Value = 0.8764025646701328
$ws.Range("Q2").Value = 9.262952336977778
$ws.Range("R2").Value = 83.3665710328
$ws.Range("S2").Value = 0.4184110992520291
$ws.Range("T2").Value = 0.418411099252029
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 5.423951666666667
$ws.Range("H3").Value = 16.271855
$ws.Range("I3").Value = 0.4774188439413272
$ws.Range("J3").Value = 0.4774188439413271
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.240846
$ws.Range("N3").Value = 0.7225379999999999
$ws.Range("O3").Value = 0.1235974353298672
$ws.Range("P3").Value = 0.1235974353298672
$ws.Range("Q3").Value = 1.30633706311
$ws.Range("R3").Value = 11.75703356799
$ws.Range("S3").Value = 0.05900774468929815
$ws.Range("T3").Value = 0.05900774468929813
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.583504333333333
$ws.Range("H4").Value = 4.750513
$ws.Range("I4").Value = 0.1393808158066948
$ws.Range("J4").Value = 0.1393808158066948
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.707786666666667
$ws.Range("N4").Value = 5.12336
$ws.Range("O4").Value = 0.8764025646701329
$ws.Range("P4").Value = 0.8764025646701328
$ws.Range("Q4").Value = 2.704287587075555
$ws.Range("R4").Value = 24.33858828368
$ws.Range("S4").Value = 0.1221537044388027
$ws.Range("T4").Value = 0.1221537044388027
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.583504333333333
$ws.Range("H5").Value = 4.750513
$ws.Range("I5").Value = 0.1393808158066948
$ws.Range("J5").Value = 0.1393808158066948
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.240846
$ws.Range("N5").Value = 0.7225379999999999
$ws.Range("O5").Value = 0.1235974353298672
$ws.Range("P5").Value = 0.1235974353298672
$ws.Range("Q5").Value = 0.3813806846659999
$ws.Range("R5").Value = 3.432426161994
$ws.Range("S5").Value = 0.0172271113678921
$ws.Range("T5").Value = 0.01722711136789209
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.6660723333333333
$ws.Range("H6").Value = 1.998217
$ws.Range("I6").Value = 0.05862800830537802
$ws.Range("J6").Value = 0.05862800830537802
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.707786666666667
$ws.Range("N6").Value = 5.12336
$ws.Range("O6").Value = 0.8764025646701329
$ws.Range("P6").Value = 0.8764025646701328
$ws.Range("Q6").Value = 1.137509449902222
$ws.Range("R6").Value = 10.23758504912
$ws.Range("S6").Value = 0.05138173684033515
$ws.Range("T6").Value = 0.05138173684033515
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.6660723333333333
$ws.Range("H7").Value = 1.998217
$ws.Range("I7").Value = 0.05862800830537802
$ws.Range("J7").Value = 0.05862800830537802
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.240846
$ws.Range("N7").Value = 0.7225379999999999
$ws.Range("O7").Value = 0.1235974353298672
$ws.Range("P7").Value = 0.1235974353298672
$ws.Range("Q7").Value = 0.160420857194
$ws.Range("R7").Value = 1.443787714746
$ws.Range("S7").Value = 0.007246271465042877
$ws.Range("T7").Value = 0.007246271465042876
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.687463666666666
$ws.Range("H8").Value = 11.062391
$ws.Range("I8").Value = 0.3245723319466
$ws.Range("J8").Value = 0.3245723319466
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.707786666666667
$ws.Range("N8").Value = 5.12336
$ws.Range("O8").Value = 0.8764025646701329
$ws.Range("P8").Value = 0.8764025646701328
$ws.Range("Q8").Value = 6.29740128375111
$ws.Range("R8").Value = 56.67661155376
$ws.Range("S8").Value = 0.2844560241389659
$ws.Range("T8").Value = 0.2844560241389659
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.687463666666666
$ws.Range("H9").Value = 11.062391
$ws.Range("I9").Value = 0.3245723319466
$ws.Range("J9").Value = 0.3245723319466
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.240846
$ws.Range("N9").Value = 0.7225379999999999
$ws.Range("O9").Value = 0.1235974353298672
$ws.Range("P9").Value = 0.1235974353298672
$ws.Range("Q9").Value = 0.8881108742619999
$ws.Range("R9").Value = 7.992997868357999
$ws.Range("S9").Value = 0.04011630780763407
$ws.Range("T9").Value = 0.04011630780763407